# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    current "2022-Q2" sheet) and populate it with the Q3 fund-holding table.
# 2. Update the "总计" (summary) sheet so the new Q3 row is on top and the
#    older rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet, positioned before "2022-Q2"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Copy the Q2 sheet's layout/formatting (header + 6 data rows share the
# same visual style as the other quarter sheets) onto the new sheet.
$q2.Range("A1:H3").Copy()
$q3.Range("A1:H3").PasteSpecial(-4122)
$q2.Range("A2:H2").Copy()
$q3.Range("A4:H7").PasteSpecial(-4122)

# ---- header row ----
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# ---- index column (A) + rank column (H) are real numbers ----
$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1
$q3.Range("A4").Value = 2
$q3.Range("A5").Value = 3
$q3.Range("A6").Value = 4
$q3.Range("A7").Value = 5

$q3.Range("H2").Value = 8
$q3.Range("H3").Value = 2
$q3.Range("H4").Value = 8
$q3.Range("H5").Value = 2
$q3.Range("H6").Value = 9
$q3.Range("H7").Value = 9

# ---- text columns (B,C,D,E,F,G) must stay text (fund codes keep leading
# zeros, decimals must render exactly like the source data) -> stage the
# values in a scratch 6x6 block formatted as Text, then copy just the
# values (not the formatting) onto the real cells so the destination keeps
# its original (unstyled) look. The scratch block must have the exact same
# shape (6 rows x 6 cols) as the destination or Excel will mis-tile values.
$scratch = $q3.Range("Z1:AE6")
$scratch.NumberFormat = "@"

$scratch.Cells.Item(1,1).Value = "001543"
$scratch.Cells.Item(1,2).Value = "宝盈新锐灵活配置混合A"
$scratch.Cells.Item(1,3).Value = "2.20"
$scratch.Cells.Item(1,4).Value = "91.14"
$scratch.Cells.Item(1,5).Value = "4.44"
$scratch.Cells.Item(1,6).Value = "0.0977"

$scratch.Cells.Item(2,1).Value = "006323"
$scratch.Cells.Item(2,2).Value = "合煦智远嘉选混合A"
$scratch.Cells.Item(2,3).Value = "0.78"
$scratch.Cells.Item(2,4).Value = "73.07"
$scratch.Cells.Item(2,5).Value = "5.99"
$scratch.Cells.Item(2,6).Value = "0.0467"

$scratch.Cells.Item(3,1).Value = "007578"
$scratch.Cells.Item(3,2).Value = "宝盈新锐灵活配置混合C"
$scratch.Cells.Item(3,3).Value = "0.36"
$scratch.Cells.Item(3,4).Value = "91.14"
$scratch.Cells.Item(3,5).Value = "4.44"
$scratch.Cells.Item(3,6).Value = "0.0160"

$scratch.Cells.Item(4,1).Value = "006324"
$scratch.Cells.Item(4,2).Value = "合煦智远嘉选混合C"
$scratch.Cells.Item(4,3).Value = "0.14"
$scratch.Cells.Item(4,4).Value = "73.07"
$scratch.Cells.Item(4,5).Value = "5.99"
$scratch.Cells.Item(4,6).Value = "0.0084"

$scratch.Cells.Item(5,1).Value = "007288"
$scratch.Cells.Item(5,2).Value = "合煦智远消费主题股票C"
$scratch.Cells.Item(5,3).Value = "0.47"
$scratch.Cells.Item(5,4).Value = "20.37"
$scratch.Cells.Item(5,5).Value = "0.88"
$scratch.Cells.Item(5,6).Value = "0.0041"

$scratch.Cells.Item(6,1).Value = "007287"
$scratch.Cells.Item(6,2).Value = "合煦智远消费主题股票A"
$scratch.Cells.Item(6,3).Value = "0.05"
$scratch.Cells.Item(6,4).Value = "20.37"
$scratch.Cells.Item(6,5).Value = "0.88"
$scratch.Cells.Item(6,6).Value = "0.0004"

$scratch.Copy()
$q3.Range("B2:G7").PasteSpecial(-4163)
$scratch.Clear()

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet -- insert the new Q3 row at the top of
# the data (row 2) and push the existing rows down.
# ---------------------------------------------------------------------
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 5
$summary.Range("D5").Value = 1.16

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 4
$summary.Range("D4").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.03

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.17

# row 5's index cell needs the same style as the other index cells (A2:A4)
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)
$summary.Range("A5").Value = 3
